# Add 2022-Q3 data
# 1) Insert a new "2022-Q3" worksheet right after "2022-Q2" -> wait, it should go BEFORE "2022-Q2"
#    (position 2, right after "总计"), carrying the new quarter's fund holdings.
# 2) Prepend a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" (totals) sheet - shift existing data rows down by one and
# insert the new 2022-Q3 summary row at the top of the data (row 2).
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Shift rows 2:7 -> 3:8 (copies values + styles together, keeping the A-column
# index style intact on every row).
$totals.Range("A2:D7").Copy($totals.Range("A3:D8"))

# Write the new first data row (2022-Q3).
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.03

# Recompute the running index in column A for the rows that shifted down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
$totals.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# Part 2: new "2022-Q3" fund-holdings sheet, inserted right after "总计"
# (i.e. before the existing "2022-Q2" sheet), with the same layout/style as
# the other quarterly sheets.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Clone layout + formatting from the "2022-Q2" sheet, then overwrite values.
$q2.Range("A1:H3").Copy($q3.Range("A1:H3"))

# Header row (unchanged labels, already copied) - data rows below.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'011686"
$q3.Range("C2").Value = "创金合信先进装备股票C"
$q3.Range("D2").Value = "'0.57"
$q3.Range("E2").Value = "'80.17"
$q3.Range("F2").Value = "'3.38"
$q3.Range("G2").Value = "'0.0193"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'011685"
$q3.Range("C3").Value = "创金合信先进装备股票A"
$q3.Range("D3").Value = "'0.25"
$q3.Range("E3").Value = "'80.17"
$q3.Range("F3").Value = "'3.38"
$q3.Range("G3").Value = "'0.0084"
$q3.Range("H3").Value = 10

Write-Output "2022-Q3 sheet added and totals updated"
